$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "respiration" sheet (sheet3): add the dark-bottle note in I2, and update
#    the F2 selection later (after other sheets exist) so window state ends
#    up correct. Adding this string first means it gets shared-string index
#    35, ahead of the o2_change_mg/L + app (mg/hr*L) renames below.
# ---------------------------------------------------------------------------
$wsResp = $wb.Worksheets.Item("respiration")
$wsResp.Range("I2").Value = "same thing as respiration, but select dark bottles instead"

# ---------------------------------------------------------------------------
# 2. "data" sheet (sheet1): rename O1/P1 headers (o2_change_mg -> o2_change_mg/L,
#    app (mg/hr) -> app (mg/hr*L)). N1 keeps its text (time_difference_hr) and
#    is automatically re-indexed once the unused strings are garbage collected.
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")
$wsData.Range("O1").Value = "o2_change_mg/L"
$wsData.Range("P1").Value = "app (mg/hr*L)"

# ---------------------------------------------------------------------------
# 3. "anpp" sheet (sheet2): fill in the first analysis row (inflow / light
#    bottle averages) plus the placeholder sample_event_id column below it.
# ---------------------------------------------------------------------------
$wsAnpp = $wb.Worksheets.Item("anpp")
$wsAnpp.Range("A2").Value = 1
$wsAnpp.Range("B2").Value = "TR_FRW1"
$wsAnpp.Range("C2").Value = "inflow"
$wsAnpp.Range("D2").Value = 42877
$wsAnpp.Range("D2").NumberFormat = "m/d/yy"
$wsAnpp.Range("E2").Value = "light"
$wsAnpp.Range("F2").Formula = "=AVERAGE(data!P2:P7)"
$wsAnpp.Range("F2").NumberFormat = "0.00"
$wsAnpp.Range("G2").Formula = "=(STDEV(data!P2:P7))/SQRT(COUNT(data!P2:P7))"
$wsAnpp.Range("A3").Value = 2
$wsAnpp.Range("A4").Value = 3
$wsAnpp.Range("A5").Value = 4
$wsAnpp.Range("A6").Value = 5
$wsAnpp.Range("A7").Value = 6
$wsAnpp.Range("A8").Value = 7

# ---------------------------------------------------------------------------
# 4. Add the two new sheets "location" and "seasons" after "respiration",
#    each with the same sample_event_id/site/location/date/bottle-type header
#    row used by the anpp/respiration sheets.
# ---------------------------------------------------------------------------
$wsLocation = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsLocation.Name = "location"
$wsLocation.Range("A1").Value = "sample_event_id"
$wsLocation.Range("B1").Value = "site"
$wsLocation.Range("C1").Value = "location"
$wsLocation.Range("D1").Value = "date"
$wsLocation.Range("E1").Value = "bottle type"
$wsLocation.PageSetup.LeftMargin = 54
$wsLocation.PageSetup.RightMargin = 54
$wsLocation.PageSetup.TopMargin = 72
$wsLocation.PageSetup.BottomMargin = 72
$wsLocation.PageSetup.HeaderMargin = 36
$wsLocation.PageSetup.FooterMargin = 36

$wsSeasons = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsSeasons.Name = "seasons"
$wsSeasons.Range("A1").Value = "sample_event_id"
$wsSeasons.Range("B1").Value = "site"
$wsSeasons.Range("C1").Value = "location"
$wsSeasons.Range("D1").Value = "date"
$wsSeasons.Range("E1").Value = "bottle type"
$wsSeasons.PageSetup.LeftMargin = 54
$wsSeasons.PageSetup.RightMargin = 54
$wsSeasons.PageSetup.TopMargin = 72
$wsSeasons.PageSetup.BottomMargin = 72
$wsSeasons.PageSetup.HeaderMargin = 36
$wsSeasons.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 5. Restore per-sheet selections / view state to match the edited workbook.
# ---------------------------------------------------------------------------
$wsData.Activate()
$wsData.Range("D38:D45").Select()

$wsAnpp.Activate()
$wsAnpp.Range("A1:E1").Select()
$excel.ActiveWindow.Zoom = 150

$wsResp.Activate()
$wsResp.Range("F2").Select()

$wsLocation.Activate()
$wsLocation.Range("A1:E1").Select()

$wsSeasons.Activate()
$wsSeasons.Range("A1:E1").Select()

$excel.ActiveWindow.Left = 2740
$excel.ActiveWindow.Top = 1020
$excel.ActiveWindow.Width = 34740
$excel.ActiveWindow.Height = 22560
